# Updated symbol list on Thu Feb  9 07:59:20 UTC 2023 with GitHub Actions
# Refresh cryptocurrency Price (D) and Volume(1h) (E) columns with latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "322.05"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-3.01%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "42.87"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-6.45%"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-7.48%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08194"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.90%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.325"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-2.92%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.828"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-10.43%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9354"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-3.87%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1112"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-4.40%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1863"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-2.98%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09469"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-4.84%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04619"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.08%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.420"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-28.41%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.38%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001289"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.46%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-5.63%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.361"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.35%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.43%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3337"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.91%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1388"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.16%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2623"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.11%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04161"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.92%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001247"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-4.91%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004309"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-6.03%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-15.58%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0002978"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-20.52%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02723"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-1.56%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05559"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-4.30%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007969"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.77%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1392"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-3.07%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006550"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-10.22%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002092"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "3.87%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007498"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-6.76%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3530"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "3.84%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006985"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-4.10%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.22%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003471"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-0.96%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003529"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.69%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.22%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001999"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.22%"
